$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.28
$ws.Range("G2").Value = 1.31
$ws.Range("H2").Value = 13
$ws.Range("J2").Value = 6
$ws.Range("Q2").Value = 1.61
$ws.Range("U2").Value = 1.75
$ws.Range("W2").Value = 4.2
$ws.Range("AF2").Value = 8.4
